$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 294, shifting existing rows 294-405 down to 295-406.
$ws.Rows(294).Insert()

# Populate the newly inserted row 294 with the new record.
$ws.Cells.Item(294, 1).Value = 4
$ws.Cells.Item(294, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(294, 3).Value = "Los Lagos"
$ws.Cells.Item(294, 4).Value = 44825
$ws.Cells.Item(294, 5).Value = 10
$ws.Cells.Item(294, 6).Value = 100112023
$ws.Cells.Item(294, 7).Value = "Brócoli"
$ws.Cells.Item(294, 8).Value = "Sin especificar"
$ws.Cells.Item(294, 9).Value = "Primera"
$ws.Cells.Item(294, 10).Value = 500
$ws.Cells.Item(294, 11).Value = 1500
$ws.Cells.Item(294, 12).Value = 1500
$ws.Cells.Item(294, 13).Value = 1500
$ws.Cells.Item(294, 14).Value = "$/unidad"
$ws.Cells.Item(294, 15).Value = "Región Metropolitana"
$ws.Cells.Item(294, 16).Value = 1500
$ws.Cells.Item(294, 17).Value = 1
$ws.Cells.Item(294, 18).Value = "Hortaliza"
